$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.850.30"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.76%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.937.42"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.90%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "243.59"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4899"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.99%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06892"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  +0.60%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "105.18"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.57%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07786"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.935.46"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.92%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.350"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.88%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.7006"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "274.99"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.26%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "30.844.43"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.86%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000007715"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.70%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.07"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "2.193.57"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.57%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.572"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.33%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.08%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.539"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.37%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.868"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "165.72"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.13%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "19.61"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.36%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.161"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.39%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.1043"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.31%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.391"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.30%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.555"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.83%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.567"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.382"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.62%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.04888"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.64%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.7606"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.42%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.150"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.61%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.9990"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.733"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.16%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.02006"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  -1.78%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.517"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.90%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "78.35"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +5.71%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.091"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.88%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.9082"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.03%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.4440"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.38%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "107.92"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.30%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.04%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.699"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -5.82%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.004.40"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.40%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.1250"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.37%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "36.25"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.56%  "
